$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text: "02-12-2018" -> "02/12/18"
#    Touches the Slide Master, every Custom Layout, and the Notes Master.
# ---------------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "02-12-2018") {
                $tr.Characters(1, $tr.Length).Text = "02/12/18"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape($master.Shapes)

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DateShape($layout.Shapes)
}

$notesMaster = $p.NotesMaster
Update-DateShape($notesMaster.Shapes)

# ---------------------------------------------------------------------------
# 2) Slide 3 content updates (JDoc-comment wording pass).
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)

# "Content Placeholder 3": only the first paragraph's text changes, the
# other two paragraphs ("Serialization-Deserialization using threading" and
# "Firebase linkup and companion android app") stay untouched.
$shpContent3 = $slide3.Shapes.Item("Content Placeholder 3")
$trContent3 = $shpContent3.TextFrame.TextRange
$oldFirstPara = "Main menu, Leaderboards, Store, Resume game and Profile and their functionalities."
$trContent3.Characters(1, $oldFirstPara.Length).Text = "All Static GUI components."

# "Content Placeholder 5": the single paragraph becomes three paragraphs.
$shpContent5 = $slide3.Shapes.Item("Content Placeholder 5")
$trContent5 = $shpContent5.TextFrame.TextRange
$trContent5.Text = "All animation related components. `rAll logic for gameplay.`rSerialization-Deserialization"
